# Update Name of Algo
# Applies updated KNN-imputed values to the result_data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "B2" = 5.965
    "A3" = -21.62
    "E3" = 16.231
    "E12" = 17.896
    "A14" = -21.5
    "A16" = -21.696
    "B18" = 5.427
    "A21" = -20.325
    "A23" = -20.536
    "B24" = 4.833
    "E24" = 17.019
    "A25" = -21.217
    "B25" = 5.854
    "E25" = 17.149
    "A26" = -21.268
    "B27" = 5.615
    "A29" = -21.116
    "B30" = 6.410000000000001
    "B31" = 6.047000000000001
    "B39" = 7.179
    "A40" = -20.391
    "E41" = 16.526
    "B42" = 8.242999999999999
    "B48" = 5.274
    "E50" = 16.442
    "B51" = 6.802
    "B52" = 5.952
    "A53" = -21.257
    "E53" = 16.581
    "B55" = 5.593
    "B56" = 5.5
    "E56" = 16.236
    "A57" = -22.053
    "B57" = 5.415999999999999
    "E57" = 16.416
    "E58" = 16.569
    "A59" = -22.404
    "B60" = 5.831
    "E61" = 16.854
    "E63" = 17.661
    "E64" = 17.642
    "A65" = -21.418
    "A69" = -21.44
    "E70" = 17.661
    "E72" = 17.027
    "B73" = 6.804
    "B74" = 8.565000000000001
    "A79" = -21.251
    "A83" = -21.338
    "E86" = 16.638
    "B89" = 4.942
    "E89" = 17.328
    "B90" = 5.767
    "A91" = -21.462
    "B92" = 5.853
    "A93" = -21.32399999999999
    "E98" = 16.244
    "A100" = -21.672
    "E100" = 16.514
    "E102" = 16.446
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

Write-Output ("Updated {0} cells in {1}" -f $updates.Count, $ws.Name)
